$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Platform")

# Rename test at row 10 (was "GeneratingAQuoteForALeadTest")
$ws.Range("A10").Value = "VerifyGenerateQuoteBackButtonTest"

# Rows 12-14 are no longer the last "Y" rows, flip their Status to "N"
$ws.Range("D12").Value = "N"
$ws.Range("D13").Value = "N"
$ws.Range("D14").Value = "N"

# Add 5 new rows (15-19), copying row 14's formatting/style so the new
# rows keep the same look (bordered/aligned "Oneplus" style column C, etc.)
$ws.Rows.Item(14).Copy()
$ws.Rows.Item(15).Insert()
$ws.Rows.Item(14).Copy()
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(14).Copy()
$ws.Rows.Item(17).Insert()
$ws.Rows.Item(14).Copy()
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(14).Copy()
$ws.Rows.Item(19).Insert()

# Fill the new rows' test case names (columns B/C already "Android"/"Oneplus"
# via the row copy); Status stays "N" for all but the final new row.
$ws.Range("A15").Value = "GenerateQuoteSearchFunctionalityTest"
$ws.Range("B15").Value = "Android"
$ws.Range("C15").Value = "Oneplus"
$ws.Range("D15").Value = "N"

$ws.Range("A16").Value = "SelectingALeadFromGenerateQuotePageTest"
$ws.Range("B16").Value = "Android"
$ws.Range("C16").Value = "Oneplus"
$ws.Range("D16").Value = "N"

$ws.Range("A17").Value = "GenerateQuoteCloseValidationTest"
$ws.Range("B17").Value = "Android"
$ws.Range("C17").Value = "Oneplus"
$ws.Range("D17").Value = "N"

$ws.Range("A18").Value = "GenerateQuoteContinueButtonCloseValidationTest"
$ws.Range("B18").Value = "Android"
$ws.Range("C18").Value = "Oneplus"
$ws.Range("D18").Value = "N"

$ws.Range("A19").Value = "GenerateQuoteExitButtonValidationTest"
$ws.Range("B19").Value = "Android"
$ws.Range("C19").Value = "Oneplus"
$ws.Range("D19").Value = "Y"

# Match the saved selection/active cell from the authored edit
$ws.Activate()
$ws.Range("D18").Select()
